$p = $ppt.ActivePresentation

# 1. Slide 16: table's style id changes from the default
#    "{19D290EE-0118-4270-B004-FB7AFBD87B77}" table style to
#    "{BDF8CD4B-A115-4F44-A5B6-7D6D636D39B8}" (a different built-in
#    PowerPoint table style). Table styles cannot be assigned through
#    the Style property - ApplyStyle must be used instead.
$slide16 = $p.Slides.Item(16)
$tableShape = $slide16.Shapes.Item(3)
$table = $tableShape.Table
$table.ApplyStyle("{BDF8CD4B-A115-4F44-A5B6-7D6D636D39B8}")

# 2. The deck's theme color scheme switches from the custom "Integral"
#    palette to the stock "Office" palette (the colours that used to
#    live only on the Notes Master's theme). Update every theme colour
#    slot through the presentation's theme colour scheme.
$themeColors = $p.Slides.Item(1).ThemeColorScheme
$themeColors.Item(1).RGB  = RGB(0x00, 0x00, 0x00)   ' dk1
$themeColors.Item(2).RGB  = RGB(0xFF, 0xFF, 0xFF)   ' lt1
$themeColors.Item(3).RGB  = RGB(0x44, 0x54, 0x6A)   ' dk2
$themeColors.Item(4).RGB  = RGB(0xE7, 0xE6, 0xE6)   ' lt2
$themeColors.Item(5).RGB  = RGB(0x5B, 0x9B, 0xD5)   ' accent1
$themeColors.Item(6).RGB  = RGB(0xED, 0x7D, 0x31)   ' accent2
$themeColors.Item(7).RGB  = RGB(0xA5, 0xA5, 0xA5)   ' accent3
$themeColors.Item(8).RGB  = RGB(0xFF, 0xC0, 0x00)   ' accent4
$themeColors.Item(9).RGB  = RGB(0x44, 0x72, 0xC4)   ' accent5
$themeColors.Item(10).RGB = RGB(0x70, 0xAD, 0x47)   ' accent6
$themeColors.Item(11).RGB = RGB(0x05, 0x63, 0xC1)   ' hlink
$themeColors.Item(12).RGB = RGB(0x95, 0x4F, 0x72)   ' folHlink
